$d = $word.ActiveDocument

$replacements = @(
    @{old = "975÷2="; new = "804÷8="},
    @{old = "779÷7="; new = "388÷7="},
    @{old = "662÷7="; new = "196÷8="},
    @{old = "909÷2="; new = "908÷6="},
    @{old = "117÷7="; new = "162÷8="},
    @{old = "176÷3="; new = "854÷3="},
    @{old = "454÷2="; new = "277÷4="},
    @{old = "966÷9="; new = "387÷3="},
    @{old = "859÷6="; new = "878÷7="},
    @{old = "166÷9="; new = "444÷3="},
    @{old = "534÷6="; new = "503÷6="},
    @{old = "606÷6="; new = "127÷3="},
    @{old = "534÷7="; new = "923÷8="},
    @{old = "791÷9="; new = "120÷8="},
    @{old = "776÷9="; new = "842÷6="},
    @{old = "885÷5="; new = "767÷2="},
    @{old = "346÷3="; new = "710÷3="},
    @{old = "651÷2="; new = "539÷9="},
    @{old = "343÷3="; new = "477÷6="},
    @{old = "259÷3="; new = "430÷3="},
    @{old = "112÷6="; new = "839÷4="},
    @{old = "411÷3="; new = "892÷8="},
    @{old = "719÷8="; new = "343÷9="},
    @{old = "266÷8="; new = "398÷7="},
    @{old = "533÷5="; new = "519÷7="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
